# Add a new "massa" layer row (row 13) to the "camadas" worksheet,
# matching the style/formatting already used by the other data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting (bold/border style on column A, blank cells in
# C/D) of the last existing data row (row 12) into the new row 13.
$ws.Range("A12:D12").Copy($ws.Range("A13:D13"))

# Now fill in the actual values for the new "massa" layer.
$ws.Range("A13").Value = "massa"
$ws.Range("B13").Value = "geometry, cd_identificador_hidrografia_poligono, nm_acidente, cd_tipo_acidente"
